# Blind model names in the "Evaluations" sheet, column C (rows 2-49),
# replacing real model identifiers with blinded labels Model A-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

$map = @{
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro"    = "Model B"
    "gpt-5.1"         = "Model C"
    "kimi-k2"         = "Model D"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
